# Add "Bazzano Davide" team as a new row (row 55) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = "Bazzano Davide"
$ws.Range("B55").Value = "Alberto Cerisara | Shark Attack"
$ws.Range("C55").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("D55").Value = "Federico Andreis | iMontagna"
$ws.Range("E55").Value = "Leonardo Trinco | Vigili del Fusto"
$ws.Range("F55").Value = "Daniele Ruzzenenti | Demobusters"
